# REPORT_20XX.X-BRXX.docx — "Alterando visualização de entrada de dados"
#
# 1) Remove the leading "CAIXA1, CAIXA2, " run pair from the intro paragraph.
# 2) Replace "*Nome da caixa*" (intro paragraph) with two bold runs:
#    "CELULA 1, " and "CELULA 2".
# 3) Rename the three "ENTIDADE N DA CAIXA M" bullet items to "...CELULA M".
# 4) Rename the three "CAIXAx - ENTIDADE N DA CAIXA x;" bullet items to
#    "CELULA x - ENTIDADE N DA CELULA x; " (note trailing space added).
# 5) In the R2 question, split "*Nome da caixa* é de ([-13.54, -13.44, -17.39])"
#    into "CELULA 1, " / "CELULA 2" (bold) + " é de " (non-bold), dropping the
#    bracketed numbers entirely.

$d = $word.ActiveDocument

# --- 1) Drop "CAIXA1, " + "CAIXA2, " runs entirely -------------------------
$d.Content.Find.Execute("CAIXA1, CAIXA2, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# --- 2) "*Nome da caixa*" -> "CELULA 1, " + "CELULA 2" (two bold runs) -----
$rng = $d.Content
$rng.Find.Execute("*Nome da caixa*", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$start = $rng.Start
$rng.Text = "CELULA 2"

$rngBefore = $d.Range($start, $start)
$rngBefore.InsertBefore("CELULA 1, ")
$rngNew = $d.Range($start, $start + 10)
$rngNew.Bold = $true

# --- 3) Bullet list: ENTIDADE N DA CAIXA M -> ENTIDADE N DA CELULA M -------
$d.Content.Find.Execute("ENTIDADE 1 DA CAIXA 1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ENTIDADE 1 DA CELULA 1", 2) | Out-Null
$d.Content.Find.Execute("ENTIDADE 2 DA CAIXA 1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ENTIDADE 2 DA CELULA 1", 2) | Out-Null
$d.Content.Find.Execute("ENTIDADE 1 DA CAIXA 2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ENTIDADE 1 DA CELULA 2", 2) | Out-Null

# --- 4) Bullet list: "CAIXAx - ENTIDADE N DA CAIXA x;" -> "CELULA x - ... ; "
$d.Content.Find.Execute("CAIXA1 - ENTIDADE 1 DA CAIXA 1;", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CELULA 1 - ENTIDADE 1 DA CELULA 1; ", 2) | Out-Null
$d.Content.Find.Execute("CAIXA1 - ENTIDADE 2 DA CAIXA 1;", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CELULA 1 - ENTIDADE 2 DA CELULA 1; ", 2) | Out-Null
$d.Content.Find.Execute("CAIXA2 - ENTIDADE 1 DA CAIXA 2;", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CELULA 2 - ENTIDADE 1 DA CELULA 2; ", 2) | Out-Null

# --- 5) R2 question: split the caixa-name + drop the bracketed numbers -----
$rng2 = $d.Content
$rng2.Find.Execute("*Nome da caixa* é de ([-13.54, -13.44, -17.39]) dBm", `
                    $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$start2 = $rng2.Start
$rng2.Text = "dBm"

$rngTail = $d.Range($start2, $start2)
$rngTail.InsertBefore(" é de ")
$rngTail2 = $d.Range($start2, $start2 + 6)
$rngTail2.Bold = $false

$rngBefore2 = $d.Range($start2, $start2)
$rngBefore2.InsertBefore("CELULA 2")
$rngCelula2 = $d.Range($start2, $start2 + 8)
$rngCelula2.Bold = $true

$rngBefore1 = $d.Range($start2, $start2)
$rngBefore1.InsertBefore("CELULA 1, ")
$rngCelula1 = $d.Range($start2, $start2 + 10)
$rngCelula1.Bold = $true

Write-Output "done"
